$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the header row (row 1), shifting all data rows up by one.
$ws.Rows.Item(1).Delete()

# Update the selection / view to match the post-edit state.
$ws.Range("J14").Select()
